$d = $word.ActiveDocument

# --- Step 1: merge the split "closure" sentence into a single run and
# drop the stray bookmark that used to sit mid-sentence (it gets
# re-added, anchored at the very end of the document, in step 3).
# The paragraph's *displayed* text already reads as the final sentence
# (runs are just split around the bookmark), so re-assigning the same
# string is a no-op; go through a throwaway placeholder first to force
# the runs to actually collapse into one. ---
$find = $d.Content.Find
$find.Execute("I cued for my closure as I went over my impact slide. T") | Out-Null
$closurePara = $find.Parent.Paragraphs(1)
$pr = $closurePara.Range
$closureStart = $pr.Start
$pr2 = $d.Range($closureStart, $pr.End - 1)
$pr2.Text = "TEMP_PLACEHOLDER_45EC4A78"
$finalClosureText = "I cued for my closure as I went over my impact slide. That was the transition or cue to my conclusion slide."
$pr3 = $d.Range($closureStart, $closureStart + "TEMP_PLACEHOLDER_45EC4A78".Length)
$pr3.Text = $finalClosureText

# --- Step 2: insert the new "mispronounced words" paragraph, right
# after "My areas of strength..." + the blank line that follows it. ---
$find2 = $d.Content.Find
$find2.Execute("My areas of strength is that I did not read off the slide") | Out-Null
$strengthPara = $find2.Parent.Paragraphs(1)
$blankPara = $strengthPara.Next()
$blankPara.Range.InsertParagraphAfter()
$newPara = $blankPara.Next()
$newPara.Range.InsertAfter("I mispronounced some words and I am coughing due to allergies here. I can do better by practicing certain words.")

# --- Step 3: append the new paragraph with the YouTube link, carrying
# the bookmark that used to live in the "closure" sentence. A trailing
# placeholder character is used so the bookmark can be anchored right
# after the link text (collapsing a range at the literal end of the
# document's last paragraph does not stick), then the placeholder is
# deleted, leaving a bookmark that wraps zero characters at the end. ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$linkPara = $d.Paragraphs($d.Paragraphs.Count)
$linkPara.Range.InsertAfter("https://youtu.be/EtV-eBcDRy8Z")
$linkPara = $d.Paragraphs($d.Paragraphs.Count)
$fullLinkRange = $linkPara.Range
$placeholderRange = $d.Range($fullLinkRange.End - 2, $fullLinkRange.End - 1)
$d.Bookmarks.Add("_GoBack", $placeholderRange) | Out-Null
$placeholderRange2 = $d.Range($fullLinkRange.End - 2, $fullLinkRange.End - 1)
$placeholderRange2.Text = ""

Write-Host "Done"
